# Autofit header columns and bold the header row on every sheet.
#
# "Grade" sheet (sheet1): headers in A1:C1 -> Class Name / Student Info / Grade
#   target column widths: A=12, B=14, C=7
# All other sheets (Algebra/Trigonometry/Geometry/Calculus/Statistics):
#   headers in A1:D1 -> Last Name / First Name / ID Number / Grade
#   target column widths: A=11, B=12, C=11, D=7
#
# Note: Excel's stored (OOXML) column width is the value assigned to
# ColumnWidth plus 5/6 (the default character-padding offset for the
# Calibri 11 "Normal" style). To land on an exact integer stored width we
# back that offset out before assigning ColumnWidth.
$offset = 5 / 6

function Set-ExactColumnWidth($ws, $colLetter, $targetWidth) {
    $ws.Columns($colLetter).ColumnWidth = $targetWidth - $offset
}

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name

    if ($name -eq "Grades") {
        $headerRange = $ws.Range("A1:C1")
        $widths = @{ "A" = 12; "B" = 14; "C" = 7 }
    } else {
        $headerRange = $ws.Range("A1:D1")
        $widths = @{ "A" = 11; "B" = 12; "C" = 11; "D" = 7 }
    }

    # Bold the header row.
    $headerRange.Font.Bold = $true

    # Autofit (exact) column widths for the header columns.
    foreach ($col in $widths.Keys) {
        Set-ExactColumnWidth $ws $col $widths[$col]
    }
}
